$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-02-11 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-12 Thursday", 2) | Out-Null
$d.Content.Find.Execute("92-24=68", $true, $false, $false, $false, $false, $true, 1, $false, "28+24=52", 2) | Out-Null
$d.Content.Find.Execute("50-15=35", $true, $false, $false, $false, $false, $true, 1, $false, "58+24=82", 2) | Out-Null
$d.Content.Find.Execute("60-46=14", $true, $false, $false, $false, $false, $true, 1, $false, "29+63=92", 2) | Out-Null
$d.Content.Find.Execute("53-26=27", $true, $false, $false, $false, $false, $true, 1, $false, "37+29=66", 2) | Out-Null
$d.Content.Find.Execute("92-67=25", $true, $false, $false, $false, $false, $true, 1, $false, "29+5=34", 2) | Out-Null
$d.Content.Find.Execute("29+67=96", $true, $false, $false, $false, $false, $true, 1, $false, "96-89=7", 2) | Out-Null
$d.Content.Find.Execute("9+49=58", $true, $false, $false, $false, $false, $true, 1, $false, "80-67=13", 2) | Out-Null
$d.Content.Find.Execute("47+16=63", $true, $false, $false, $false, $false, $true, 1, $false, "78-29=49", 2) | Out-Null
$d.Content.Find.Execute("93-37=56", $true, $false, $false, $false, $false, $true, 1, $false, "60-38=22", 2) | Out-Null
$d.Content.Find.Execute("27+38=65", $true, $false, $false, $false, $false, $true, 1, $false, "37+57=94", 2) | Out-Null
$d.Content.Find.Execute("48+44=92", $true, $false, $false, $false, $false, $true, 1, $false, "50-31=19", 2) | Out-Null
$d.Content.Find.Execute("36-27=9", $true, $false, $false, $false, $false, $true, 1, $false, "66-9=57", 2) | Out-Null
$d.Content.Find.Execute("40-14=26", $true, $false, $false, $false, $false, $true, 1, $false, "85-58=27", 2) | Out-Null
$d.Content.Find.Execute("25+66=91", $true, $false, $false, $false, $false, $true, 1, $false, "2+19=21", 2) | Out-Null
$d.Content.Find.Execute("43-37=6", $true, $false, $false, $false, $false, $true, 1, $false, "56+19=75", 2) | Out-Null
$d.Content.Find.Execute("40-12=28", $true, $false, $false, $false, $false, $true, 1, $false, "63-35=28", 2) | Out-Null
$d.Content.Find.Execute("70-66=4", $true, $false, $false, $false, $false, $true, 1, $false, "22+19=41", 2) | Out-Null
$d.Content.Find.Execute("68+27=95", $true, $false, $false, $false, $false, $true, 1, $false, "5+66=71", 2) | Out-Null
$d.Content.Find.Execute("86-37=49", $true, $false, $false, $false, $false, $true, 1, $false, "49+39=88", 2) | Out-Null
$d.Content.Find.Execute("19+59=78", $true, $false, $false, $false, $false, $true, 1, $false, "68+26=94", 2) | Out-Null
$d.Content.Find.Execute("33+49=82", $true, $false, $false, $false, $false, $true, 1, $false, "61-54=7", 2) | Out-Null
$d.Content.Find.Execute("96-57=39", $true, $false, $false, $false, $false, $true, 1, $false, "30-14=16", 2) | Out-Null
$d.Content.Find.Execute("6+27=33", $true, $false, $false, $false, $false, $true, 1, $false, "90-53=37", 2) | Out-Null
$d.Content.Find.Execute("40-4=36", $true, $false, $false, $false, $false, $true, 1, $false, "47+49=96", 2) | Out-Null
$d.Content.Find.Execute("36+7=43", $true, $false, $false, $false, $false, $true, 1, $false, "91-33=58", 2) | Out-Null
$d.Content.Find.Execute("90-19=71", $true, $false, $false, $false, $false, $true, 1, $false, "95-29=66", 2) | Out-Null
$d.Content.Find.Execute("13+79=92", $true, $false, $false, $false, $false, $true, 1, $false, "90-3=87", 2) | Out-Null
$d.Content.Find.Execute("80-12=68", $true, $false, $false, $false, $false, $true, 1, $false, "7+17=24", 2) | Out-Null
$d.Content.Find.Execute("59+38=97", $true, $false, $false, $false, $false, $true, 1, $false, "30-5=25", 2) | Out-Null
$d.Content.Find.Execute("9+32=41", $true, $false, $false, $false, $false, $true, 1, $false, "22+39=61", 2) | Out-Null
$d.Content.Find.Execute("70-11=59", $true, $false, $false, $false, $false, $true, 1, $false, "51-24=27", 2) | Out-Null
$d.Content.Find.Execute("23+69=92", $true, $false, $false, $false, $false, $true, 1, $false, "28+15=43", 2) | Out-Null
$d.Content.Find.Execute("90-36=54", $true, $false, $false, $false, $false, $true, 1, $false, "4+69=73", 2) | Out-Null
$d.Content.Find.Execute("77+6=83", $true, $false, $false, $false, $false, $true, 1, $false, "69+24=93", 2) | Out-Null
$d.Content.Find.Execute("18+74=92", $true, $false, $false, $false, $false, $true, 1, $false, "90-2=88", 2) | Out-Null
$d.Content.Find.Execute("39+44=83", $true, $false, $false, $false, $false, $true, 1, $false, "96-87=9", 2) | Out-Null
$d.Content.Find.Execute("57-18=39", $true, $false, $false, $false, $false, $true, 1, $false, "39+57=96", 2) | Out-Null
$d.Content.Find.Execute("98-19=79", $true, $false, $false, $false, $false, $true, 1, $false, "47+5=52", 2) | Out-Null
$d.Content.Find.Execute("8+3=11", $true, $false, $false, $false, $false, $true, 1, $false, "18+14=32", 2) | Out-Null
$d.Content.Find.Execute("24+48=72", $true, $false, $false, $false, $false, $true, 1, $false, "25+18=43", 2) | Out-Null
$d.Content.Find.Execute("69+14=83", $true, $false, $false, $false, $false, $true, 1, $false, "58-29=29", 2) | Out-Null
$d.Content.Find.Execute("92-35=57", $true, $false, $false, $false, $false, $true, 1, $false, "43-16=27", 2) | Out-Null
$d.Content.Find.Execute("39+19=58", $true, $false, $false, $false, $false, $true, 1, $false, "72-39=33", 2) | Out-Null
$d.Content.Find.Execute("8+89=97", $true, $false, $false, $false, $false, $true, 1, $false, "17+39=56", 2) | Out-Null
$d.Content.Find.Execute("65+28=93", $true, $false, $false, $false, $false, $true, 1, $false, "82-29=53", 2) | Out-Null
$d.Content.Find.Execute("55+8=63", $true, $false, $false, $false, $false, $true, 1, $false, "61-12=49", 2) | Out-Null
$d.Content.Find.Execute("26+56=82", $true, $false, $false, $false, $false, $true, 1, $false, "96-19=77", 2) | Out-Null
$d.Content.Find.Execute("28+57=85", $true, $false, $false, $false, $false, $true, 1, $false, "4+79=83", 2) | Out-Null
$d.Content.Find.Execute("81-25=56", $true, $false, $false, $false, $false, $true, 1, $false, "36+57=93", 2) | Out-Null
$d.Content.Find.Execute("18+9=27", $true, $false, $false, $false, $false, $true, 1, $false, "45+18=63", 2) | Out-Null
$d.Content.Find.Execute("45+38=83", $true, $false, $false, $false, $false, $true, 1, $false, "50-5=45", 2) | Out-Null
$d.Content.Find.Execute("17+47=64", $true, $false, $false, $false, $false, $true, 1, $false, "36+19=55", 2) | Out-Null
$d.Content.Find.Execute("19+62=81", $true, $false, $false, $false, $false, $true, 1, $false, "19+43=62", 2) | Out-Null
$d.Content.Find.Execute("33+39=72", $true, $false, $false, $false, $false, $true, 1, $false, "17+28=45", 2) | Out-Null
$d.Content.Find.Execute("33-28=5", $true, $false, $false, $false, $false, $true, 1, $false, "93-7=86", 2) | Out-Null
$d.Content.Find.Execute("79+5=84", $true, $false, $false, $false, $false, $true, 1, $false, "61-43=18", 2) | Out-Null
$d.Content.Find.Execute("6+6=12", $true, $false, $false, $false, $false, $true, 1, $false, "74-16=58", 2) | Out-Null
$d.Content.Find.Execute("92-14=78", $true, $false, $false, $false, $false, $true, 1, $false, "35+46=81", 2) | Out-Null
$d.Content.Find.Execute("81-66=15", $true, $false, $false, $false, $false, $true, 1, $false, "90-73=17", 2) | Out-Null
$d.Content.Find.Execute("49+33=82", $true, $false, $false, $false, $false, $true, 1, $false, "19+54=73", 2) | Out-Null
$d.Content.Find.Execute("71-44=27", $true, $false, $false, $false, $false, $true, 1, $false, "17+36=53", 2) | Out-Null
$d.Content.Find.Execute("86-58=28", $true, $false, $false, $false, $false, $true, 1, $false, "39+32=71", 2) | Out-Null
$d.Content.Find.Execute("71-56=15", $true, $false, $false, $false, $false, $true, 1, $false, "13+58=71", 2) | Out-Null
$d.Content.Find.Execute("42-25=17", $true, $false, $false, $false, $false, $true, 1, $false, "76-58=18", 2) | Out-Null
$d.Content.Find.Execute("51-43=8", $true, $false, $false, $false, $false, $true, 1, $false, "70-53=17", 2) | Out-Null
$d.Content.Find.Execute("48+14=62", $true, $false, $false, $false, $false, $true, 1, $false, "40-32=8", 2) | Out-Null
$d.Content.Find.Execute("9+72=81", $true, $false, $false, $false, $false, $true, 1, $false, "91-38=53", 2) | Out-Null
$d.Content.Find.Execute("71-55=16", $true, $false, $false, $false, $false, $true, 1, $false, "54-5=49", 2) | Out-Null
$d.Content.Find.Execute("87+6=93", $true, $false, $false, $false, $false, $true, 1, $false, "35-8=27", 2) | Out-Null
$d.Content.Find.Execute("52-35=17", $true, $false, $false, $false, $false, $true, 1, $false, "6+29=35", 2) | Out-Null
$d.Content.Find.Execute("29+58=87", $true, $false, $false, $false, $false, $true, 1, $false, "45+29=74", 2) | Out-Null
$d.Content.Find.Execute("49+48=97", $true, $false, $false, $false, $false, $true, 1, $false, "91-63=28", 2) | Out-Null
$d.Content.Find.Execute("79+3=82", $true, $false, $false, $false, $false, $true, 1, $false, "58-29=29", 2) | Out-Null
$d.Content.Find.Execute("23-4=19", $true, $false, $false, $false, $false, $true, 1, $false, "41-16=25", 2) | Out-Null
$d.Content.Find.Execute("85-46=39", $true, $false, $false, $false, $false, $true, 1, $false, "42-35=7", 2) | Out-Null
$d.Content.Find.Execute("83-29=54", $true, $false, $false, $false, $false, $true, 1, $false, "27+37=64", 2) | Out-Null
$d.Content.Find.Execute("48+27=75", $true, $false, $false, $false, $false, $true, 1, $false, "57-9=48", 2) | Out-Null
$d.Content.Find.Execute("27+35=62", $true, $false, $false, $false, $false, $true, 1, $false, "17+6=23", 2) | Out-Null
$d.Content.Find.Execute("88+9=97", $true, $false, $false, $false, $false, $true, 1, $false, "92-49=43", 2) | Out-Null
$d.Content.Find.Execute("39+6=45", $true, $false, $false, $false, $false, $true, 1, $false, "70-5=65", 2) | Out-Null
$d.Content.Find.Execute("82-69=13", $true, $false, $false, $false, $false, $true, 1, $false, "38+14=52", 2) | Out-Null
$d.Content.Find.Execute("26-19=7", $true, $false, $false, $false, $false, $true, 1, $false, "50-7=43", 2) | Out-Null
$d.Content.Find.Execute("39+4=43", $true, $false, $false, $false, $false, $true, 1, $false, "7+69=76", 2) | Out-Null
$d.Content.Find.Execute("63-59=4", $true, $false, $false, $false, $false, $true, 1, $false, "73+18=91", 2) | Out-Null
$d.Content.Find.Execute("76-49=27", $true, $false, $false, $false, $false, $true, 1, $false, "83-75=8", 2) | Out-Null
$d.Content.Find.Execute("7+27=34", $true, $false, $false, $false, $false, $true, 1, $false, "83-56=27", 2) | Out-Null
$d.Content.Find.Execute("19+65=84", $true, $false, $false, $false, $false, $true, 1, $false, "57+15=72", 2) | Out-Null
$d.Content.Find.Execute("83-78=5", $true, $false, $false, $false, $false, $true, 1, $false, "85-36=49", 2) | Out-Null
$d.Content.Find.Execute("19+75=94", $true, $false, $false, $false, $false, $true, 1, $false, "39+14=53", 2) | Out-Null
$d.Content.Find.Execute("16+46=62", $true, $false, $false, $false, $false, $true, 1, $false, "23+18=41", 2) | Out-Null
$d.Content.Find.Execute("7+59=66", $true, $false, $false, $false, $false, $true, 1, $false, "80-54=26", 2) | Out-Null
$d.Content.Find.Execute("74-49=25", $true, $false, $false, $false, $false, $true, 1, $false, "72+9=81", 2) | Out-Null
$d.Content.Find.Execute("55+17=72", $true, $false, $false, $false, $false, $true, 1, $false, "59+16=75", 2) | Out-Null
$d.Content.Find.Execute("68-29=39", $true, $false, $false, $false, $false, $true, 1, $false, "90-84=6", 2) | Out-Null
$d.Content.Find.Execute("39+59=98", $true, $false, $false, $false, $false, $true, 1, $false, "9+74=83", 2) | Out-Null
$d.Content.Find.Execute("31-16=15", $true, $false, $false, $false, $false, $true, 1, $false, "18+5=23", 2) | Out-Null
$d.Content.Find.Execute("71-46=25", $true, $false, $false, $false, $false, $true, 1, $false, "70-57=13", 2) | Out-Null
$d.Content.Find.Execute("67-8=59", $true, $false, $false, $false, $false, $true, 1, $false, "58+4=62", 2) | Out-Null
$d.Content.Find.Execute("53+19=72", $true, $false, $false, $false, $false, $true, 1, $false, "35+8=43", 2) | Out-Null
